# Weekly update: a new price observation was added at the top of the
# data (row 4), pushing the existing rows 4-7 down to rows 5-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 4-7 down by one row, inserting a blank row at 4.
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new weekly observation.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44791
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 8500
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8750
$ws.Range("N4").Value = "$/cuna 10 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 875
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"
